$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grant_multi_vouchers")

$ws.Range("B12").Value = "\用户红包券\\用户红包券\\用户红包券\\用户红包券\\用户红包券\\用户红包券\\用户红包券\"

$ws.Range("C8:C12").Formula = "=RANDBETWEEN(1,10)"
$ws.Range("D8:D12").Formula = "=RANDBETWEEN(1,9998)"

$ws.Range("E8").Formula = "=RANDBETWEEN(D8,9999)"
$ws.Range("E9").Formula = "=RANDBETWEEN(D9,9999)"
$ws.Range("E10").Formula = "=RANDBETWEEN(D10,9999)"
$ws.Range("E11").Formula = "=RANDBETWEEN(D11,9999)"
$ws.Range("E12").Formula = "=RANDBETWEEN(D12,9999)"

$ws.Range("C12").Select()
